$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("groups")

# Convert the string "false" values in D2:D6 into real boolean FALSE values
$ws.Range("D2:D6").Value = $false
